$wb = $excel.ActiveWorkbook

$oldGuid = "1cb91d82-c797-4adc-bdaf-aeab64b649ef"
$newGuid = "b4402a84-1312-432e-9eb6-26acbbf4d091"
$oldHash = "f90b542ae74e50539d3d5afef1203f0252100bff"
$newHash = "d4a67a5bffdab32b6531ae11dee13103bc2cfe8d"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet: B2 is the "Path And Name" hyperlink text ---
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# --- zh-cn sheet ---
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-22 15:09:36"
$wsZh.Range("I2").ClearContents()
$wsZh.Range("J2").ClearContents()
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# --- de-de sheet ---
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-22 15:09:41"
$wsDe.Range("I2").ClearContents()
$wsDe.Range("J2").ClearContents()
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

# Update remaining hyperlink display text (A2 hyperlinks) on zh-cn/de-de sheets
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

# Autofit columns I and J on zh-cn/de-de after clearing content
$wsZh.Columns.Item(9).AutoFit() | Out-Null
$wsZh.Columns.Item(10).AutoFit() | Out-Null
$wsDe.Columns.Item(9).AutoFit() | Out-Null
$wsDe.Columns.Item(10).AutoFit() | Out-Null
